# Auto-generated edit script applying cryptos list update (2024-06-28)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.266.47"
$ws.Range("E2").Value = "  -1.91%  "
$ws.Range("D3").Value = "'3.377.73"
$ws.Range("E3").Value = "  -2.09%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'571.15"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("D6").Value = "'140.44"
$ws.Range("E6").Value = "  -5.16%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'3.377.81"
$ws.Range("E8").Value = "  -2.16%  "
$ws.Range("D9").Value = "'0.474"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "'7.46"
$ws.Range("E10").Value = "  -4.46%  "
$ws.Range("D11").Value = "'0.123"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").Value = "'0.393"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "'3.957.85"
$ws.Range("E13").Value = "  -2.02%  "
$ws.Range("D14").Value = "'28.13"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "'0.0000170"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("D17").Value = "'3.375.74"
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("D18").Value = "'60.498.80"
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("D19").Value = "'6.26"
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("D20").Value = "'14.08"
$ws.Range("E20").Value = "  -1.82%  "
$ws.Range("D21").Value = "'9.21"
$ws.Range("E21").Value = "  -2.66%  "
$ws.Range("D22").Value = "'387.66"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("D24").Value = "'73.45"
$ws.Range("E24").Value = "  +1.02%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "'0.0000116"
$ws.Range("E26").Value = "  -4.72%  "
$ws.Range("D27").Value = "'3.518.46"
$ws.Range("E27").Value = "  -2.15%  "
$ws.Range("E28").Value = "  -0.94%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "'7.38"
$ws.Range("E30").Value = "  -5.68%  "
$ws.Range("D31").Value = "'8.05"
$ws.Range("E31").Value = "  -2.36%  "
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("D33").Value = "'1.43"
$ws.Range("E33").Value = "  -6.23%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "'23.70"
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("D36").Value = "'6.92"
$ws.Range("E36").Value = "  -2.20%  "
$ws.Range("D37").Value = "'3.410.29"
$ws.Range("E37").Value = "  -1.79%  "
$ws.Range("D38").Value = "'167.27"
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "'1.50"
$ws.Range("E39").Value = "  -4.24%  "
$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").Value = "'4.91"
$ws.Range("E40").Value = "  -6.04%  "
$ws.Range("D41").Value = "'0.0775"
$ws.Range("E41").Value = "  -1.79%  "
$ws.Range("D42").Value = "'27.04"
$ws.Range("E42").Value = "  +3.68%  "
$ws.Range("D43").Value = "'0.781"
$ws.Range("E43").Value = "  -1.92%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").Value = "'4.45"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("D46").Value = "'1.69"
$ws.Range("E46").Value = "  -1.80%  "
$ws.Range("D47").Value = "'41.31"
$ws.Range("E47").Value = "  -2.38%  "
$ws.Range("D48").Value = "'2.530.28"
$ws.Range("E48").Value = "  -2.91%  "
$ws.Range("D49").Value = "'1.11"
$ws.Range("E49").Value = "  -4.00%  "
$ws.Range("D50").Value = "'6.83"
$ws.Range("E50").Value = "  -2.29%  "
$ws.Range("D51").Value = "'23.01"
$ws.Range("E51").Value = "  -0.96%  "
